$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 3.229285597997556
$ws.Range("R2").Value = 29.063570381978
$ws.Range("S2").Value = 0.000061659803659281023770037666
$ws.Range("T2").Value = 0.000061659803659281037322564822

$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 657.6094073941372
$ws.Range("R3").Value = 5918.484666547234
$ws.Range("S3").Value = 0.01255635827613453
$ws.Range("T3").Value = 0.01255635827613453

$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 125.65356481492
$ws.Range("R4").Value = 1130.88208333428
$ws.Range("S4").Value = 0.00239922233585689
$ws.Range("T4").Value = 0.00239922233585689

$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 200.361154482186
$ws.Range("R5").Value = 1803.250390339673
$ws.Range("S5").Value = 0.003825684991745289
$ws.Range("T5").Value = 0.003825684991745289

$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("S6").Value = 0.7790597529863222
$ws.Range("T6").Value = 0.7790597529863222

$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 7796.180469987221
$ws.Range("R7").Value = 70165.62422988498
$ws.Range("S7").Value = 0.1488598460816897
$ws.Range("T7").Value = 0.1488598460816897

$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 11.44810952125011
$ws.Range("R8").Value = 103.032985691251
$ws.Range("S8").Value = 0.0002185895808620771
$ws.Range("T8").Value = 0.0002185895808620771

$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 2331.284827430789
$ws.Range("R9").Value = 20981.5634468771
$ws.Range("S9").Value = 0.04451342576276901
$ws.Range("T9").Value = 0.04451342576276901

$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 445.45325214614
$ws.Range("R10").Value = 4009.07926931526
$ws.Range("S10").Value = 0.008505460180961048
$ws.Range("T10").Value = 0.008505460180961048

